$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.226.59'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.860.28'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.33'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2870'
$ws.Range("E8").Value = '  +0.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06549'
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.77'
$ws.Range("E10").Value = '  +3.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07928'
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.71'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.870.28'
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.181'
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6820'
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '268.38'
$ws.Range("E16").Value = '  -5.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.223.75'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.78'
$ws.Range("E18").Value = '  +8.68%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007443'
$ws.Range("E20").Value = '  +2.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.115.70'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.330'
$ws.Range("E22").Value = '  -3.48%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.189'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.12'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.229'
$ws.Range("E26").Value = '  -0.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.92'
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.961'
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.385'
$ws.Range("E29").Value = '  +2.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09851'
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.390'
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.474'
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04713'
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.136'
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7035'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.616'
$ws.Range("E39").Value = '  +3.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.252'
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.65'
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8464'
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4176'
$ws.Range("E44").Value = '  -0.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9992'
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.29'
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '956.83'
$ws.Range("E47").Value = '  +2.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.175'
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.183'
$ws.Range("E49").Value = '  -1.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.15'
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05663'
$ws.Range("E51").Value = '  +0.53%  '
